$wb = $excel.ActiveWorkbook

# --- 1) Rename header cells on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3) Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the bold/centered/bordered header style already used on the
# "Weekly Quantity" sheet (style index 1) instead of creating new styles.
$wsWeekly.Range("A1:D1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# --- 4) Data rows ---
$forecastData = @(
    @(44941.99999999999, 0, -36.73949017274813, 27.89219140227145),
    @(44955.99999999999, 2, -32.12695504560176, 35.46348858967492),
    @(44976.99999999999, 12, -22.51081409642472, 43.03504121366336),
    @(45011.99999999999, 27, -6.100031583922957, 60.63733130106394),
    @(45046.99999999999, 43, 11.77572139319088, 76.42377997250794),
    @(45060.99999999999, 49, 14.85372592143417, 84.19209019047004),
    @(45067.99999999999, 52, 20.29403788927094, 87.088723453797),
    @(45095.99999999999, 65, 32.99561400192061, 95.94285948229093),
    @(45123.99999999999, 77, 47.09623475683514, 112.7953656202602),
    @(45130.99999999999, 80, 48.08562691282658, 111.5928192564785),
    @(45137.99999999999, 83, 47.61101258775901, 112.5923597389446),
    @(45144.99999999999, 86, 53.63040718301654, 118.5785043723998),
    @(45158.99999999999, 93, 60.30281764216173, 125.729399540547),
    @(45165.99999999999, 96, 60.50447831244848, 131.1563366786123),
    @(45172.99999999999, 99, 66.99829969094907, 131.6828018117286),
    @(45179.99999999999, 102, 68.92935954540265, 134.3902131868925),
    @(45186.99999999999, 105, 70.17265475611909, 138.9423790203081),
    @(45193.99999999999, 108, 75.83904762228896, 142.6470496607432),
    @(45200.99999999999, 111, 78.86513188231451, 143.0591932332247),
    @(45207.99999999999, 114, 83.14554024522529, 146.8211118399818),
    @(45214.99999999999, 118, 82.99601390105904, 149.5540014434308)
)

$row = 2
foreach ($entry in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Column A on the data rows uses the same date-formatted style (index 2)
# as column A on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A22").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()
